$d = $word.ActiveDocument

function Find-Replace($searchText, $replaceText) {
    $find = $d.Content.Find
    $find.ClearFormatting()
    $find.Replacement.ClearFormatting()
    return $find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)
}

# 1) ${departmentFull} -> ${support}
#    Unlike the other placeholders below, this one keeps its original
#    three-run split ("${" / "support" / "}") instead of collapsing into a
#    single run, so replace only the inner word...
$r1 = Find-Replace "`${departmentFull}" "`${support}"
Write-Output "departmentFull->support: $r1"

# ...then re-split the merged run back into three runs (matching the
# original "${" / word / "}" run boundaries) by nudging and reverting a
# character property on the middle word only. Also drops the now-stale
# spellcheck proofErr markers around it.
$find2 = $d.Content.Find
$find2.ClearFormatting()
$found2 = $find2.Execute("`${support}", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found2) {
    $rng = $find2.Parent
    $s = $rng.Start
    $e = $rng.End
    $innerLen = ($e - $s) - 3   # exclude the leading "${" and trailing "}"
    $mid = $d.Range($s + 2, $s + 2 + $innerLen)
    $mid.Bold = 1
    $mid.Bold = 0
}

# 2) The remaining placeholders collapse their "${" / name / "}" runs
#    (and drop the spellcheck proofErr markers around the name) into a
#    single run each.
$r2 = Find-Replace "`${ilos}" "`${ilos}"
Write-Output "ilos: $r2"

$r3 = Find-Replace "`${budgetSource}" "`${budgetSource}"
Write-Output "budgetSource: $r3"

$r4 = Find-Replace "`${sig_cscp}" "`${sig_cscp}"
Write-Output "sig_cscp: $r4"

$r5 = Find-Replace "`${sig_csca}" "`${sig_csca}"
Write-Output "sig_csca: $r5"

$r6 = Find-Replace "`${sig_sscp}" "`${sig_sscp}"
Write-Output "sig_sscp: $r6"

$r7 = Find-Replace "`${sig_dean}" "`${sig_dean}"
Write-Output "sig_dean: $r7"
